{"js": "// Resize/reposition the floating \"Text Box 3\" text box (the notes box that\n// explains the 1-5 numeric rating scale) to match the author's manual\n// drag-resize: it grew taller (more line-wrap room) and shifted slightly\n// right/down.\nconst shapes = context.document.body.shapes;\nshapes.load(\"items/name\");\nawait context.sync();\n\nconst textBox = shapes.items.find((s) => s.name === \"Text Box 3\");\nif (!textBox) {\n  throw new Error('Shape \"Text Box 3\" was not found in the document body.');\n}\n\n// Target geometry, expressed in points (1 pt = 12700 EMU) so the saved\n// OOXML lands exactly on the EMU values from the edit:\n//   left   -> posOffsetH 3716655 -> 3721281 EMU\n//   top    -> posOffsetV  219075 ->  225185 EMU\n//   height -> extent cy   699770 ->  795647 EMU (width cx stays 2138680)\ntextBox.left = 3721281 / 12700;\ntextBox.top = 225185 / 12700;\ntextBox.height = 795647 / 12700;\n\nawait context.sync();\n", "ps1": "# Resize/reposition the floating \"Text Box 3\" text box (the notes box that\n# explains the 1-5 numeric rating scale) to match the author's manual\n# drag-resize: it grew taller (more line-wrap room) and shifted slightly\n# right/down.\n$d = $word.ActiveDocument\n\n$textBox = $null\nforeach ($shp in $d.Shapes) {\n    if ($shp.Name -eq \"Text Box 3\") {\n        $textBox = $shp\n        break\n    }\n}\n\nif ($null -eq $textBox) {\n    throw 'Shape \"Text Box 3\" was not found in the document body.'\n}\n\n# Target geometry, expressed in points (1 pt = 12700 EMU) so the saved\n# OOXML lands exactly on the EMU values from the edit:\n#   Left   -> posOffsetH 3716655 -> 3721281 EMU\n#   Top    -> posOffsetV  219075 ->  225185 EMU\n#   Height -> extent cy   699770 ->  795647 EMU (width cx stays 2138680)\n$textBox.Left = 3721281 / 12700\n$textBox.Top = 225185 / 12700\n$textBox.Height = 795647 / 12700\n"}
